# Add a "mount info" option: the long boilerplate paragraph explaining the
# allowable withdrawal force of the lag screw is replaced with a single
# template placeholder token "mountInfoDef", while keeping the paragraph's
# formatting (pPr/rPr) and the proofErr markers that already bracketed the
# "lbs" word intact.

$d = $word.ActiveDocument

# Locate the single, unique occurrence of "lbs" (inside the target run).
$hit = $d.Content
$found = $hit.Find.Execute("lbs", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'lbs' run to anchor the replacement."
}
$lbsStart = $hit.Start
$lbsEnd = $hit.End

# Grab the paragraph that contains it.
$para = $hit.Paragraphs(1)
$paraStart = $para.Range.Start
$paraEnd = $para.Range.End

# Remove everything in the paragraph after the "lbs" run (but keep the
# trailing paragraph mark itself, hence "- 1").
$tail = $d.Range($lbsEnd, $paraEnd - 1)
if ($tail.Start -lt $tail.End) {
    $tail.Text = ""
}

# Remove everything in the paragraph before the "lbs" run.
$head = $d.Range($paraStart, $lbsStart)
if ($head.Start -lt $head.End) {
    $head.Text = ""
}

# The remaining run (originally "lbs") now sits right at the paragraph
# start; replace its text with the new placeholder token.
$wordLen = $lbsEnd - $lbsStart
$final = $d.Range($paraStart, $paraStart + $wordLen)
$final.Text = "mountInfoDef"
